$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 1862.25
$ws.Range("J53").Value = 3438.8333
$ws.Range("L53").Value = 3438.8333
$ws.Range("N53").Value = -4712.8333
$ws.Range("H129").Value = 2408
$ws.Range("J129").Value = 2659.3333
$ws.Range("L129").Value = 7977.999899999999
$ws.Range("N129").Value = -17977.9999
$ws.Range("H131").Value = 6483.174
$ws.Range("I131").Value = 4374.4546
$ws.Range("K131").Value = 13123.3638
$ws.Range("M131").Value = -8083.363799999999
$ws.Range("H132").Value = 2193.093
$ws.Range("I132").Value = 2080.561
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 6241.683000000001
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -3711.683000000001
$ws.Range("N132").Value = -18560
$ws.Range("H138").Value = 5367.492
$ws.Range("J138").Value = 5380.5933
$ws.Range("L138").Value = 16141.7799
$ws.Range("N138").Value = -26421.7799
$ws.Range("H141").Value = 1692
$ws.Range("I141").Value = 1130.4
$ws.Range("K141").Value = 3391.2
$ws.Range("M141").Value = 1788.8

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8687.519
$ws.Range("I32").Value = 7790.615
$ws.Range("K32").Value = 7790.615
$ws.Range("M32").Value = -7503.615
$ws.Range("H33").Value = 14005.2
$ws.Range("I33").Value = 10008.667
$ws.Range("K33").Value = 10008.667
$ws.Range("M33").Value = -9679.666999999999
$ws.Range("H36").Value = 15427.125
$ws.Range("I36").Value = 7694
$ws.Range("J36").Value = 18004.834
$ws.Range("K36").Value = 7694
$ws.Range("L36").Value = 18004.834
$ws.Range("M36").Value = -7348
$ws.Range("N36").Value = -18696.834
$ws.Range("H45").Value = 4251
$ws.Range("I45").Value = 2996.6667
$ws.Range("K45").Value = 2996.6667
$ws.Range("M45").Value = -2619.6667

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 10257
$ws.Range("I82").Value = 10257
$ws.Range("K82").Value = 10257
$ws.Range("M82").Value = -9874
$ws.Range("H85").Value = 10257
$ws.Range("I85").Value = 10257
$ws.Range("K85").Value = 10257
$ws.Range("M85").Value = -8931
$ws.Range("H86").Value = 4649
$ws.Range("I86").Value = 3762.111
$ws.Range("J86").Value = 6929.5713
$ws.Range("K86").Value = 3762.111
$ws.Range("L86").Value = 6929.5713
$ws.Range("M86").Value = -2639.111
$ws.Range("N86").Value = -9175.5713
$ws.Range("H89").Value = 4649
$ws.Range("I89").Value = 3762.111
$ws.Range("J89").Value = 6929.5713
$ws.Range("K89").Value = 18810.555
$ws.Range("L89").Value = 34647.85649999999
$ws.Range("M89").Value = -13194.555
$ws.Range("N89").Value = -45879.85649999999
$ws.Range("H105").Value = 20455.941
$ws.Range("I105").Value = 26187.375
$ws.Range("K105").Value = 26187.375
$ws.Range("M105").Value = -24440.375

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 24294.143
$ws.Range("H51").Value = 39999
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").Value = ""
$ws.Range("H60").Value = 20875
$ws.Range("H61").Value = 39999
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").Value = ""
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = ""
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = ""
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").Value = ""
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").Value = ""
$ws.Range("H99").Value = 4886.875
$ws.Range("I99").Value = 4749.1665
$ws.Range("K99").Value = 4749.1665
$ws.Range("M99").Value = -3251.1665
$ws.Range("H106").Value = 34037.5
$ws.Range("I106").Value = 61000
$ws.Range("K106").Value = 61000
$ws.Range("M106").Value = -59738
$ws.Range("H108").Value = 82487
$ws.Range("J108").Value = 82487
$ws.Range("L108").Value = 82487
$ws.Range("N108").Value = -90167
$ws.Range("H126").Value = 4886.875
$ws.Range("I126").Value = 4749.1665
$ws.Range("K126").Value = 14247.4995
$ws.Range("M126").Value = -11777.4995

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2915.6191
$ws.Range("I3").Value = 901.55554
$ws.Range("J3").Value = 15000
$ws.Range("K3").Value = 2704.66662
$ws.Range("L3").Value = 45000
$ws.Range("M3").Value = -2592.66662
$ws.Range("N3").Value = -45224
$ws.Range("H20").Value = 4949
$ws.Range("I20").Value = 4949
$ws.Range("K20").Value = 14847
$ws.Range("M20").Value = -14620
$ws.Range("H33").Value = 199.625
$ws.Range("I33").Value = 251.4
$ws.Range("J33").Value = 113.333336
$ws.Range("K33").Value = 1508.4
$ws.Range("L33").Value = 680.000016
$ws.Range("M33").Value = -1225.4
$ws.Range("N33").Value = -1246.000016
$ws.Range("H107").Value = 1645488
$ws.Range("I107").Value = 535.5454999999999
$ws.Range("J107").Value = 3907297.5
$ws.Range("K107").Value = 1606.6365
$ws.Range("L107").Value = 11721892.5
$ws.Range("M107").Value = 313.3635000000002
$ws.Range("N107").Value = -11725732.5
$ws.Range("H113").Value = 2779.25
$ws.Range("I113").Value = 2080.5
$ws.Range("J113").Value = 2842.7727
$ws.Range("K113").Value = 6241.5
$ws.Range("L113").Value = 8528.3181
$ws.Range("M113").Value = -4071.5
$ws.Range("N113").Value = -12868.3181
$ws.Range("H131").Value = 14522252
$ws.Range("I131").Value = 41668456
$ws.Range("J131").Value = 10236010
$ws.Range("K131").Value = 125005368
$ws.Range("L131").Value = 30708030
$ws.Range("M131").Value = -125000328
$ws.Range("N131").Value = -30718110
$ws.Range("H132").Value = 5676.7856
$ws.Range("I132").Value = 4999.8335
$ws.Range("J132").Value = 6184.5
$ws.Range("K132").Value = 44998.5015
$ws.Range("L132").Value = 55660.5
$ws.Range("M132").Value = -42468.5015
$ws.Range("N132").Value = -60720.5
$ws.Range("H134").Value = 6764.05
$ws.Range("I134").Value = 6764.05
$ws.Range("K134").Value = 20292.15
$ws.Range("M134").Value = -15222.15

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 15003
$ws.Range("I21").Value = 15003
$ws.Range("K21").Value = 15003
$ws.Range("M21").Value = -14830
$ws.Range("H30").Value = 15003
$ws.Range("I30").Value = 15003
$ws.Range("K30").Value = 15003
$ws.Range("M30").Value = -14898
$ws.Range("H68").Value = 49981.5
$ws.Range("J68").Value = 49978.668
$ws.Range("L68").Value = 49978.668
$ws.Range("N68").Value = -51600.668
$ws.Range("H71").Value = 49981.5
$ws.Range("J71").Value = 49978.668
$ws.Range("L71").Value = 149936.004
$ws.Range("N71").Value = -158048.004
$ws.Range("H80").Value = 1257704
$ws.Range("J80").Value = 10272
$ws.Range("L80").Value = 10272
$ws.Range("N80").Value = -12268
$ws.Range("H83").Value = 1257704
$ws.Range("J83").Value = 10272
$ws.Range("L83").Value = 51360
$ws.Range("N83").Value = -61344
$ws.Range("H126").Value = 4749.4614
$ws.Range("I126").Value = 4692.5654
$ws.Range("K126").Value = 14077.6962
$ws.Range("M126").Value = -11607.6962

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10378.4
$ws.Range("I7").Value = 6902.1665
$ws.Range("J7").Value = 24283.334
$ws.Range("K7").Value = 6902.1665
$ws.Range("L7").Value = 24283.334
$ws.Range("M7").Value = -6790.1665
$ws.Range("N7").Value = -24507.334
$ws.Range("H34").Value = 5000
$ws.Range("I34").Value = 5000
$ws.Range("K34").Value = 5000
$ws.Range("M34").Value = -4828
$ws.Range("H40").Value = 10478.318
$ws.Range("I40").Value = 13214.444
$ws.Range("K40").Value = 13214.444
$ws.Range("M40").Value = -13078.444
$ws.Range("H55").Value = 7145937
$ws.Range("I55").Value = 12500140
$ws.Range("J55").Value = 7000
$ws.Range("K55").Value = 12500140
$ws.Range("L55").Value = 7000
$ws.Range("M55").Value = -12499967
$ws.Range("N55").Value = -7346
$ws.Range("H126").Value = 10378.4
$ws.Range("I126").Value = 6902.1665
$ws.Range("J126").Value = 24283.334
$ws.Range("K126").Value = 20706.4995
$ws.Range("L126").Value = 72850.00199999999
$ws.Range("M126").Value = -18236.4995
$ws.Range("N126").Value = -77790.00199999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").Value = ""
$ws.Range("H81").Value = 12215.2
$ws.Range("I81").Value = 10525
$ws.Range("K81").Value = 21050
$ws.Range("M81").Value = -19989
$ws.Range("H84").Value = 12215.2
$ws.Range("I84").Value = 10525
$ws.Range("K84").Value = 105250
$ws.Range("M84").Value = -99946
$ws.Range("H132").Value = 3066.257
$ws.Range("I132").Value = 1697.2174
$ws.Range("J132").Value = 5690.25
$ws.Range("K132").Value = 5091.6522
$ws.Range("L132").Value = 17070.75
$ws.Range("M132").Value = -2561.6522
$ws.Range("N132").Value = -22130.75
